# Updated cryptos list on Sun Feb 25 23:40:50 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    # Force the literal text into the cell (apostrophe prefix stops Excel's
    # automatic number/date inference) then restore the default "Normal"
    # style so no stray number-format/quote-prefix style is left behind.
    $ws.Range($range).Value = "'" + $text
    $ws.Range($range).Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue "D2" "51.759.94"
$ws.Range("E2").Value = "  +0.36%  "

# Row 3 - Ethereum
Set-TextValue "D3" "3.112.39"
$ws.Range("E3").Value = "  +4.08%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.05%  "

# Row 5 - BNB
Set-TextValue "D5" "388.70"
$ws.Range("E5").Value = "  +1.75%  "

# Row 6 - Solana
Set-TextValue "D6" "103.34"
$ws.Range("E6").Value = "  -0.89%  "

# Row 7 - XRP
$ws.Range("E7").Value = "  -0.46%  "

# Row 9 - Cardano
Set-TextValue "D9" "0.592"
$ws.Range("E9").Value = "  -0.94%  "

# Row 10 - Avalanche
Set-TextValue "D10" "37.41"
$ws.Range("E10").Value = "  +1.53%  "

# Row 11 - TRON
$ws.Range("E11").Value = "  -0.06%  "

# Row 12 - Dogecoin
$ws.Range("E12").Value = "  -0.04%  "

# Row 13 - WrappedliquidstakedEther2.0
Set-TextValue "D13" "3.601.28"
$ws.Range("E13").Value = "  +3.95%  "

# Row 14 - was Polkadot, now Chainlink
$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextValue "D14" "18.70"
$ws.Range("E14").Value = "  +0.92%  "

# Row 15 - was Chainlink, now Polkadot
$ws.Range("B15").Value = "Polkadot"
$ws.Range("C15").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue "D15" "7.91"
$ws.Range("E15").Value = "  +0.92%  "

# Row 16 - WrappedEther
Set-TextValue "D16" "3.089.72"
$ws.Range("E16").Value = "  +3.41%  "

# Row 17 - was Uniswap, now Polygon
$ws.Range("B17").Value = "Polygon"
$ws.Range("C17").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-TextValue "D17" "0.993"
$ws.Range("E17").Value = "  -0.34%  "

# Row 18 - was Polygon, now Uniswap
$ws.Range("B18").Value = "Uniswap"
$ws.Range("C18").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
Set-TextValue "D18" "10.99"
$ws.Range("E18").Value = "  -2.72%  "

# Row 19 - WrappedBTC
Set-TextValue "D19" "51.833.14"
$ws.Range("E19").Value = "  +0.39%  "

# Row 20 - ImmutableX
Set-TextValue "D20" "3.21"
$ws.Range("E20").Value = "  +4.21%  "

# Row 21 - InternetComputer(DFINITY)
Set-TextValue "D21" "12.47"
$ws.Range("E21").Value = "  -0.70%  "

# Row 22 - ShibaInu
Set-TextValue "D22" "0.0₃0968"
$ws.Range("E22").Value = "  +0.27%  "

# Row 23 - Litecoin
Set-TextValue "D23" "70.08"

# Row 24 - BitcoinCash
Set-TextValue "D24" "267.80"
$ws.Range("E24").Value = "  +0.08%  "

# Row 25 - PancakeSwap
$ws.Range("E25").Value = "  -3.03%  "

# Row 26 - Filecoin
Set-TextValue "D26" "8.15"
$ws.Range("E26").Value = "  +1.00%  "

# Row 27 - EthereumClassic
Set-TextValue "D27" "27.14"
$ws.Range("E27").Value = "  +3.77%  "

# Row 28 - Kaspa
$ws.Range("E28").Value = "  +0.64%  "

# Row 29 - was Dai, now RenderToken
$ws.Range("B29").Value = "RenderToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D29" "7.19"
$ws.Range("E29").Value = "  -0.24%  "

# Row 30 - was RenderToken, now Dai
$ws.Range("B30").Value = "Dai"
$ws.Range("C30").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextValue "D30" "1.00"
$ws.Range("E30").Value = "  +0.02%  "

# Row 31 - Hedera
$ws.Range("E31").Value = "  -0.56%  "

# Row 32 - Cosmos
Set-TextValue "D32" "10.37"
$ws.Range("E32").Value = "  -0.09%  "

# Row 33 - InjectiveProtocol
Set-TextValue "D33" "35.73"
$ws.Range("E33").Value = "  +2.92%  "

# Row 34 - Toncoin
$ws.Range("E34").Value = "  +0.77%  "

# Row 35 - OKB
Set-TextValue "D35" "50.29"
$ws.Range("E35").Value = "  -2.17%  "

# Row 36 - VeChain
Set-TextValue "D36" "0.0450"
$ws.Range("E36").Value = "  +1.11%  "

# Row 37 - FirstDigitalUSD
$ws.Range("E37").Value = "  -0.12%  "

# Row 38 - LidoDAOToken
Set-TextValue "D38" "3.39"
$ws.Range("E38").Value = "  +2.86%  "

# Row 39 - TheGraph
$ws.Range("E39").Value = "  +6.70%  "

# Row 40 - ARBITRUM
$ws.Range("E40").Value = "  +2.61%  "

# Row 42 - Monero
Set-TextValue "D42" "129.34"
$ws.Range("E42").Value = "  +1.59%  "

# Row 43 - Celestia
$ws.Range("E43").Value = "  -1.12%  "

# Row 44 - Stellar
$ws.Range("E44").Value = "  -0.07%  "

# Row 45 - NEARProtocol
$ws.Range("E45").Value = "  -3.77%  "

# Row 46 - EnergySwap
Set-TextValue "D46" "22.24"
$ws.Range("E46").Value = "  +3.63%  "

# Row 47 - ApeXProtocol
Set-TextValue "D47" "2.48"
$ws.Range("E47").Value = "  +5.24%  "

# Row 48 - WEMIXToken
$ws.Range("E48").Value = "  +2.23%  "

# Row 49 - Maker
Set-TextValue "D49" "2.050.33"
$ws.Range("E49").Value = "  +0.67%  "

# Row 50 - RocketPoolETH
Set-TextValue "D50" "3.416.10"
$ws.Range("E50").Value = "  +4.08%  "

# Row 51 - was Algorand, now BEAM
$ws.Range("B51").Value = "BEAM"
$ws.Range("C51").Value = "https://coinranking.com/coin/cYYMfXF4u+beam-beam"
Set-TextValue "D51" "0.0328"
$ws.Range("E51").Value = "  -0.85%  "
